$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column K header and value for the Manual DFR Omicron file test data
$ws.Range("K1").Value = "OmicronFile"
$ws.Range("K2").Value = "ManualDFRInjection.qcm"

# Update the active selection to match the saved worksheet view (K1 selected)
$ws.Range("K1").Select()
